$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2103.1667
$ws.Range("I6").Value = 2279.25
$ws.Range("K6").Value = 6837.75
$ws.Range("M6").Value = -6725.75

$ws.Range("H9").Value = 1713084.2
$ws.Range("I9").Value = 287
$ws.Range("J9").Value = 4796119.5
$ws.Range("K9").Value = 287
$ws.Range("L9").Value = 4796119.5
$ws.Range("M9").Value = -118
$ws.Range("N9").Value = -4796457.5

$ws.Range("H70").Value = 1457316.2
$ws.Range("I70").Value = 4571.4287
$ws.Range("K70").Value = 13714.2861
$ws.Range("M70").Value = -13444.2861

$ws.Range("H73").Value = 1457316.2
$ws.Range("I73").Value = 4571.4287
$ws.Range("K73").Value = 13714.2861
$ws.Range("M73").Value = -12778.2861

$ws.Range("H95").Value = 64497.75
$ws.Range("J95").Value = 64497.75
$ws.Range("L95").Value = 64497.75
$ws.Range("N95").Value = -69989.75

$ws.Range("H113").Value = 2533
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 1254

$ws.Range("H132").Value = 1355.7727
$ws.Range("J132").Value = 4500
$ws.Range("L132").Value = 13500
$ws.Range("N132").Value = -18560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 319
$ws.Range("I5").Value = 319
$ws.Range("K5").Value = 319
$ws.Range("M5").Value = -207

$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = ""

$ws.Range("H45").Value = 3232.111
$ws.Range("I45").Value = 2897.8
$ws.Range("K45").Value = 2897.8
$ws.Range("M45").Value = -2520.8

$ws.Range("H61").Value = 9018.120999999999
$ws.Range("I61").Value = 6981.5
$ws.Range("K61").Value = 6981.5
$ws.Range("M61").Value = -6769.5

$ws.Range("H74").Value = 2961.3809
$ws.Range("I74").Value = 942.9
$ws.Range("K74").Value = 942.9
$ws.Range("M74").Value = -68.89999999999998

$ws.Range("H77").Value = 2961.3809
$ws.Range("I77").Value = 942.9
$ws.Range("K77").Value = 4714.5
$ws.Range("M77").Value = -346.5

$ws.Range("H104").Value = 4435.8
$ws.Range("J104").Value = 4243.75
$ws.Range("L104").Value = 4243.75
$ws.Range("N104").Value = -11231.75

$ws.Range("H110").Value = 7354141.5
$ws.Range("I110").Value = 8621407
$ws.Range("K110").Value = 8621407
$ws.Range("M110").Value = -8619362

$ws.Range("H132").Value = 5513.9316
$ws.Range("I132").Value = 3845.125
$ws.Range("J132").Value = 9964.083000000001
$ws.Range("K132").Value = 11535.375
$ws.Range("L132").Value = 29892.249
$ws.Range("M132").Value = -9005.375
$ws.Range("N132").Value = -34952.249

$ws.Range("H136").Value = 9018.120999999999
$ws.Range("I136").Value = 6981.5
$ws.Range("K136").Value = 20944.5
$ws.Range("M136").Value = -18394.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 319
$ws.Range("I4").Value = 319
$ws.Range("K4").Value = 319
$ws.Range("M4").Value = -204

$ws.Range("H86").Value = 92070.13
$ws.Range("I86").Value = 1021.4375
$ws.Range("K86").Value = 1021.4375
$ws.Range("M86").Value = 101.5625

$ws.Range("H89").Value = 92070.13
$ws.Range("I89").Value = 1021.4375
$ws.Range("K89").Value = 5107.1875
$ws.Range("M89").Value = 508.8125

$ws.Range("H134").Value = 7508.1055
$ws.Range("I134").Value = 6090
$ws.Range("K134").Value = 18270
$ws.Range("M134").Value = -15735

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1190.8823
$ws.Range("J94").Value = 916.7857
$ws.Range("L94").Value = 916.7857
$ws.Range("N94").Value = -1818.7857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 2000
$ws.Range("J52").Value = 2000
$ws.Range("L52").Value = 6000
$ws.Range("N52").Value = -6532

$ws.Range("H61").Value = 105.11539
$ws.Range("I61").Value = 91.791664
$ws.Range("J61").Value = 265
$ws.Range("K61").Value = 275.374992
$ws.Range("L61").Value = 795
$ws.Range("M61").Value = -60.37499200000002
$ws.Range("N61").Value = -1225

$ws.Range("H62").Value = 2482.3794
$ws.Range("I62").Value = 1749.0834
$ws.Range("K62").Value = 5247.2502
$ws.Range("M62").Value = -4561.2502

$ws.Range("H65").Value = 2482.3794
$ws.Range("I65").Value = 1749.0834
$ws.Range("K65").Value = 15741.7506
$ws.Range("M65").Value = -12309.7506

$ws.Range("H109").Value = 2163.25
$ws.Range("I109").Value = 2163.25
$ws.Range("K109").Value = 6489.75
$ws.Range("M109").Value = -5449.75

$ws.Range("H116").Value = 5000
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").Value = ""

$ws.Range("H140").Value = 1017.4286
$ws.Range("I140").Value = 864.8889
$ws.Range("J140").Value = 1292
$ws.Range("K140").Value = 2594.6667
$ws.Range("L140").Value = 3876
$ws.Range("M140").Value = 2585.3333
$ws.Range("N140").Value = -14236

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1500
$ws.Range("I5").Value = 1500
$ws.Range("K5").Value = 1500
$ws.Range("M5").Value = -1388

$ws.Range("H97").Value = 1065.0625
$ws.Range("I97").Value = 999.25
$ws.Range("K97").Value = 999.25
$ws.Range("M97").Value = -503.25

$ws.Range("H132").Value = 7807.5
$ws.Range("I132").Value = 4334.3335
$ws.Range("J132").Value = 10412.375
$ws.Range("K132").Value = 13003.0005
$ws.Range("L132").Value = 31237.125
$ws.Range("M132").Value = -10473.0005
$ws.Range("N132").Value = -36297.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 225.33333
$ws.Range("I55").Value = 183.77777
$ws.Range("K55").Value = 183.77777
$ws.Range("M55").Value = -10.77777

$ws.Range("H93").Value = 733.3333
$ws.Range("I93").Value = 733.3333
$ws.Range("K93").Value = 733.3333
$ws.Range("M93").Value = 514.6667

$ws.Range("H136").Value = 4751.2886
$ws.Range("I136").Value = 4192.5264
$ws.Range("K136").Value = 12577.5792
$ws.Range("M136").Value = -10027.5792

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 10000
$ws.Range("I12").Value = 10000
$ws.Range("K12").Value = 10000
$ws.Range("M12").Value = -9858

$ws.Range("H107").Value = 2637.7778
$ws.Range("I107").Value = 1650
$ws.Range("J107").Value = 3131.6667
$ws.Range("K107").Value = 4950
$ws.Range("L107").Value = 9395.000100000001
$ws.Range("M107").Value = -3030
$ws.Range("N107").Value = -13235.0001

$ws.Range("H136").Value = 3347.75
$ws.Range("I136").Value = 1319.8948
$ws.Range("K136").Value = 3959.6844
$ws.Range("M136").Value = -1409.6844
